$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 142
    $ws.Range("F3").Value = 1716
    $ws.Range("F6").Value = 476
    $ws.Range("F9").Value = 638
}
